# Auto-generated edit script: updates currentAveragePrice / Leve profit columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("H33").Value = 1414.8125
$ws.Range("I33").Value = 262.45456
$ws.Range("K33").Value = 262.45456
$ws.Range("M33").Value = -33.45456000000001
$ws.Range("H62").Value = 7480.625
$ws.Range("I62").Value = 6118.077
$ws.Range("J62").Value = 9090.909
$ws.Range("K62").Value = 6118.077
$ws.Range("L62").Value = 9090.909
$ws.Range("M62").Value = -5494.077
$ws.Range("N62").Value = -10338.909
$ws.Range("H65").Value = 7480.625
$ws.Range("I65").Value = 6118.077
$ws.Range("J65").Value = 9090.909
$ws.Range("K65").Value = 30590.385
$ws.Range("L65").Value = 45454.545
$ws.Range("M65").Value = -27470.385
$ws.Range("N65").Value = -51694.545
$ws.Range("H116").Value = 7554.853
$ws.Range("I116").Value = 7775.6665
$ws.Range("K116").Value = 7775.6665
$ws.Range("M116").Value = -4333.6665
$ws.Range("H135").Value = 845.2941
$ws.Range("I135").Value = 864.2857
$ws.Range("J135").Value = 756.6667
$ws.Range("K135").Value = 7778.571300000001
$ws.Range("L135").Value = 6810.0003
$ws.Range("M135").Value = -5243.571300000001
$ws.Range("N135").Value = -11880.0003
$ws.Range("H138").Value = 2840.2292
$ws.Range("J138").Value = 3678.9167
$ws.Range("L138").Value = 11036.7501
$ws.Range("N138").Value = -21316.7501
$ws.Range("H141").Value = 5569.0557
$ws.Range("I141").Value = 3658.4546
$ws.Range("K141").Value = 10975.3638
$ws.Range("M141").Value = -5795.363799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 317.75
$ws.Range("I5").Value = 211
$ws.Range("K5").Value = 211
$ws.Range("M5").Value = -99
$ws.Range("H63").Value = 4724.9165
$ws.Range("I63").Value = 2099.8572
$ws.Range("K63").Value = 2099.8572
$ws.Range("M63").Value = -1413.8572
$ws.Range("H66").Value = 4724.9165
$ws.Range("I66").Value = 2099.8572
$ws.Range("K66").Value = 10499.286
$ws.Range("M66").Value = -7067.286
$ws.Range("H68").Value = 45000
$ws.Range("J68").Value = 45000
$ws.Range("L68").Value = 45000
$ws.Range("N68").Value = -46622
$ws.Range("H71").Value = 45000
$ws.Range("J71").Value = 45000
$ws.Range("L71").Value = 135000
$ws.Range("N71").Value = -143112
$ws.Range("H74").Value = 13895489
$ws.Range("I74").Value = 15879558
$ws.Range("K74").Value = 15879558
$ws.Range("M74").Value = -15878684
$ws.Range("H77").Value = 13895489
$ws.Range("I77").Value = 15879558
$ws.Range("K77").Value = 79397790
$ws.Range("M77").Value = -79393422
$ws.Range("H97").Value = 646.8461
$ws.Range("I97").Value = 795.9
$ws.Range("K97").Value = 795.9
$ws.Range("M97").Value = -299.9
$ws.Range("H132").Value = 6817.643
$ws.Range("I132").Value = 3712.3333
$ws.Range("K132").Value = 11136.9999
$ws.Range("M132").Value = -8606.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 317.75
$ws.Range("I4").Value = 211
$ws.Range("K4").Value = 211
$ws.Range("M4").Value = -96
$ws.Range("H76").Value = 32500
$ws.Range("J76").Value = 32500
$ws.Range("L76").Value = 32500
$ws.Range("N76").Value = -33130
$ws.Range("H79").Value = 32500
$ws.Range("J79").Value = 32500
$ws.Range("L79").Value = 32500
$ws.Range("N79").Value = -34684
$ws.Range("H107").Value = 2449.611
$ws.Range("I107").Value = 1998
$ws.Range("K107").Value = 1998
$ws.Range("M107").Value = -78

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 184.17392
$ws.Range("I11").Value = 181.8
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 181.8
$ws.Range("L11").Value = 200
$ws.Range("M11").Value = -41.80000000000001
$ws.Range("N11").Value = -480
$ws.Range("H17").Value = 79999.5
$ws.Range("I17").Value = 79999.5
$ws.Range("K17").Value = 79999.5
$ws.Range("M17").Value = -79825.5
$ws.Range("H18").Value = 70000
$ws.Range("J18").Value = 70000
$ws.Range("L18").Value = 70000
$ws.Range("N18").Value = -70460
$ws.Range("H22").Value = 1018.5769
$ws.Range("J22").Value = 1478.8334
$ws.Range("L22").Value = 1478.8334
$ws.Range("N22").Value = -2178.8334
$ws.Range("H86").Value = 13276.223
$ws.Range("I86").Value = 9333
$ws.Range("J86").Value = 15247.833
$ws.Range("K86").Value = 9333
$ws.Range("L86").Value = 15247.833
$ws.Range("M86").Value = -8210
$ws.Range("N86").Value = -17493.833
$ws.Range("H89").Value = 13276.223
$ws.Range("I89").Value = 9333
$ws.Range("J89").Value = 15247.833
$ws.Range("K89").Value = 46665
$ws.Range("L89").Value = 76239.16500000001
$ws.Range("M89").Value = -41049
$ws.Range("N89").Value = -87471.16500000001
$ws.Range("H105").Value = 3191.2727
$ws.Range("I105").Value = 1824.25
$ws.Range("J105").Value = 6836.6665
$ws.Range("K105").Value = 1824.25
$ws.Range("L105").Value = 6836.6665
$ws.Range("M105").Value = -77.25
$ws.Range("N105").Value = -10330.6665
$ws.Range("H132").Value = 5169.32
$ws.Range("I132").Value = 4458.0625
$ws.Range("K132").Value = 13374.1875
$ws.Range("M132").Value = -10844.1875
$ws.Range("H134").Value = 7645.393
$ws.Range("I134").Value = 3946.8333
$ws.Range("J134").Value = 14302.8
$ws.Range("K134").Value = 11840.4999
$ws.Range("L134").Value = 42908.39999999999
$ws.Range("M134").Value = -9305.499899999999
$ws.Range("N134").Value = -47978.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 17709.857
$ws.Range("I70").Value = 9950
$ws.Range("J70").Value = 19003.166
$ws.Range("K70").Value = 29850
$ws.Range("L70").Value = 57009.49800000001
$ws.Range("M70").Value = -29535
$ws.Range("N70").Value = -57639.49800000001
$ws.Range("H73").Value = 17709.857
$ws.Range("I73").Value = 9950
$ws.Range("J73").Value = 19003.166
$ws.Range("K73").Value = 29850
$ws.Range("L73").Value = 57009.49800000001
$ws.Range("M73").Value = -28758
$ws.Range("N73").Value = -59193.49800000001
$ws.Range("H129").Value = 5557617.5
$ws.Range("I129").Value = 815.36365
$ws.Range("K129").Value = 2446.09095
$ws.Range("M129").Value = 2553.90905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 44593.332
$ws.Range("I20").Value = 34005
$ws.Range("J20").Value = 47618.57
$ws.Range("K20").Value = 34005
$ws.Range("L20").Value = 47618.57
$ws.Range("M20").Value = -33760
$ws.Range("N20").Value = -48108.57
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 8507058
$ws.Range("I31").Value = 392.33334
$ws.Range("J31").Value = 11342613
$ws.Range("K31").Value = 392.33334
$ws.Range("L31").Value = 11342613
$ws.Range("M31").Value = -144.33334
$ws.Range("N31").Value = -11343109
$ws.Range("H40").Value = 8083.636
$ws.Range("I40").Value = 6151.2666
$ws.Range("J40").Value = 12224.429
$ws.Range("K40").Value = 6151.2666
$ws.Range("L40").Value = 12224.429
$ws.Range("M40").Value = -6015.2666
$ws.Range("N40").Value = -12496.429
$ws.Range("H93").Value = 12325.833
$ws.Range("I93").Value = 2326
$ws.Range("J93").Value = 22325.666
$ws.Range("K93").Value = 2326
$ws.Range("L93").Value = 22325.666
$ws.Range("M93").Value = -1078
$ws.Range("N93").Value = -24821.666
$ws.Range("H122").Value = 14502
$ws.Range("I122").Value = 9999
$ws.Range("K122").Value = 29997
$ws.Range("M122").Value = -27547
$ws.Range("H132").Value = 3990.1
$ws.Range("I132").Value = 2321.7778
$ws.Range("J132").Value = 19005
$ws.Range("K132").Value = 6965.3334
$ws.Range("L132").Value = 57015
$ws.Range("M132").Value = -4435.3334
$ws.Range("N132").Value = -62075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1686.2858
$ws.Range("I96").Value = 1400
$ws.Range("J96").Value = 1800.8
$ws.Range("K96").Value = 1400
$ws.Range("L96").Value = 1800.8
$ws.Range("M96").Value = -27
$ws.Range("N96").Value = -4546.8
$ws.Range("H132").Value = 5853.885
$ws.Range("I132").Value = 5487.15
$ws.Range("K132").Value = 16461.45
$ws.Range("M132").Value = -13931.45

